$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data (refreshed / re-sorted by Utility desc), one row shorter than before.
$data = @(
    @("property_type","beds","sum",0.78431034642846653),
    @("property_type","price","max",0.70970053391506394),
    @("property_type","number_of_reviews","avg",0.643122570016608),
    @("zipcode","number_of_reviews","sum",0.63298052157375517),
    @("property_type","number_of_reviews","sum",0.61782642262966037),
    @("zipcode","price","sum",0.53605271217625983),
    @("zipcode","price","max",0.50089181041350084),
    @("zipcode","beds","sum",0.49859673891477319),
    @("property_type","beds","max",0.45869395897345389),
    @("property_type","number_of_reviews","max",0.45202178881851229),
    @("zipcode","price","avg",0.37989148788184163),
    @("zipcode","number_of_reviews","max",0.3785129005789607),
    @("property_type","price","avg",0.37163389191361251),
    @("property_type","beds","avg",0.37118607110242507),
    @("zipcode","beds","max",0.37078707340725542),
    @("room_type","beds","max",0.36823815361030338),
    @("zipcode","number_of_reviews","avg",0.35402518884714651),
    @("room_type","number_of_reviews","avg",0.35284052801055138),
    @("room_type","price","max",0.33274609772311331),
    @("room_type","number_of_reviews","sum",0.31963329764213949),
    @("zipcode","beds","avg",0.30486563603373817),
    @("room_type","number_of_reviews","max",0.29062149156787909),
    @("room_type","price","avg",0.17690706673230089),
    @("room_type","beds","avg",0.1656988070359518),
    @("room_type","beds","sum",0.032711168813403142),
    @("room_type","price","sum",0.026904067135758991)
)

# Remove the old last row (28) entirely so the table shrinks by one row.
$ws.Rows.Item(28).Delete()

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Column A no longer has an explicit custom width / best-fit setting.
$ws.Columns.Item(1).ColumnWidth = 8.43

# Selection moves to K18 as recorded by the author.
$ws.Range("K18").Select()
